$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2180
$ws.Range("I18").Value = 1240
$ws.Range("K18").Value = 1240
$ws.Range("M18").Value = -956
$ws.Range("H38").Value = 1516.9286
$ws.Range("I38").Value = 1597.4615
$ws.Range("J38").Value = 470
$ws.Range("K38").Value = 4792.3845
$ws.Range("L38").Value = 1410
$ws.Range("M38").Value = -4420.3845
$ws.Range("N38").Value = -2154
$ws.Range("H42").Value = 4983.2856
$ws.Range("I42").Value = 856
$ws.Range("J42").Value = 8078.75
$ws.Range("K42").Value = 2568
$ws.Range("L42").Value = 24236.25
$ws.Range("M42").Value = -2338
$ws.Range("N42").Value = -24696.25
$ws.Range("H53").Value = 301.7647
$ws.Range("I53").Value = 213
$ws.Range("K53").Value = 213
$ws.Range("M53").Value = 424
$ws.Range("H55").Value = 627.75
$ws.Range("I55").Value = 703.6667
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 703.6667
$ws.Range("L55").Value = 400
$ws.Range("M55").Value = -489.6667
$ws.Range("N55").Value = -828
$ws.Range("H92").Value = 500
$ws.Range("I92").Value = 500
$ws.Range("K92").Value = 500
$ws.Range("M92").Value = 748
$ws.Range("H111").Value = 8617.6
$ws.Range("I111").Value = 8522.333000000001
$ws.Range("K111").Value = 25566.999
$ws.Range("M111").Value = -22499.999
$ws.Range("H112").Value = 2461.7144
$ws.Range("I112").Value = 1239.75
$ws.Range("K112").Value = 3719.25
$ws.Range("M112").Value = -2611.25
$ws.Range("H131").Value = 18054.666
$ws.Range("I131").Value = 18054.666
$ws.Range("K131").Value = 54163.99800000001
$ws.Range("M131").Value = -49123.99800000001
$ws.Range("H138").Value = 3408.58
$ws.Range("I138").Value = 2013.2142
$ws.Range("J138").Value = 3951.2222
$ws.Range("K138").Value = 6039.642599999999
$ws.Range("L138").Value = 11853.6666
$ws.Range("M138").Value = -899.6425999999992
$ws.Range("N138").Value = -22133.6666

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1012.4
$ws.Range("J74").Value = 1099.5
$ws.Range("L74").Value = 1099.5
$ws.Range("N74").Value = -2847.5
$ws.Range("H77").Value = 1012.4
$ws.Range("J77").Value = 1099.5
$ws.Range("L77").Value = 5497.5
$ws.Range("N77").Value = -14233.5
$ws.Range("H97").Value = 1131.2667
$ws.Range("I97").Value = 1121.25
$ws.Range("K97").Value = 1121.25
$ws.Range("M97").Value = -625.25
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H42").Value = 175000
$ws.Range("J42").Value = 175000
$ws.Range("L42").Value = 175000
$ws.Range("N42").Value = -175656
$ws.Range("H86").Value = 6733.1665
$ws.Range("I86").Value = 7801
$ws.Range("J86").Value = 1394
$ws.Range("K86").Value = 7801
$ws.Range("L86").Value = 1394
$ws.Range("M86").Value = -6678
$ws.Range("N86").Value = -3640
$ws.Range("H89").Value = 6733.1665
$ws.Range("I89").Value = 7801
$ws.Range("J89").Value = 1394
$ws.Range("K89").Value = 39005
$ws.Range("L89").Value = 6970
$ws.Range("M89").Value = -33389
$ws.Range("N89").Value = -18202
$ws.Range("H94").Value = 2310.4707
$ws.Range("I94").Value = 2268.3845
$ws.Range("K94").Value = 2268.3845
$ws.Range("M94").Value = -1817.3845
$ws.Range("H126").Value = 49999
$ws.Range("J126").Value = 49999
$ws.Range("L126").Value = 49999
$ws.Range("N126").Value = -59879

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1513.25
$ws.Range("I22").Value = 249
$ws.Range("K22").Value = 249
$ws.Range("M22").Value = 101
$ws.Range("H58").Value = 2490.5557
$ws.Range("I58").Value = 1829.75
$ws.Range("K58").Value = 1829.75
$ws.Range("M58").Value = -1626.75
$ws.Range("H132").Value = 5630.6875
$ws.Range("I132").Value = 5872.7334
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 17618.2002
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -15088.2002
$ws.Range("N132").Value = -11060
$ws.Range("H134").Value = 2062.7778
$ws.Range("I134").Value = 2062.7778
$ws.Range("K134").Value = 6188.3334
$ws.Range("M134").Value = -3653.3334
$ws.Range("H136").Value = 2490.5557
$ws.Range("I136").Value = 1829.75
$ws.Range("K136").Value = 5489.25
$ws.Range("M136").Value = -2939.25

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 54.375
$ws.Range("J12").Value = 67.40000000000001
$ws.Range("L12").Value = 202.2
$ws.Range("N12").Value = -548.2
$ws.Range("H68").Value = 4137.7417
$ws.Range("I68").Value = 3999
$ws.Range("J68").Value = 4142.3667
$ws.Range("K68").Value = 11997
$ws.Range("L68").Value = 12427.1001
$ws.Range("M68").Value = -11186
$ws.Range("N68").Value = -14049.1001
$ws.Range("H71").Value = 4137.7417
$ws.Range("I71").Value = 3999
$ws.Range("J71").Value = 4142.3667
$ws.Range("K71").Value = 35991
$ws.Range("L71").Value = 37281.3003
$ws.Range("M71").Value = -31935
$ws.Range("N71").Value = -45393.3003
$ws.Range("H86").Value = 239.8
$ws.Range("I86").Value = 200
$ws.Range("J86").Value = 266.33334
$ws.Range("K86").Value = 600
$ws.Range("L86").Value = 799.0000200000001
$ws.Range("M86").Value = 586
$ws.Range("N86").Value = -3171.00002
$ws.Range("H89").Value = 239.8
$ws.Range("I89").Value = 200
$ws.Range("J89").Value = 266.33334
$ws.Range("K89").Value = 1800
$ws.Range("L89").Value = 2397.00006
$ws.Range("M89").Value = 4128
$ws.Range("N89").Value = -14253.00006
$ws.Range("H113").Value = 1619.3334
$ws.Range("J113").Value = 1619.3334
$ws.Range("L113").Value = 4858.0002
$ws.Range("N113").Value = -9198.0002
$ws.Range("H129").Value = 1252682
$ws.Range("I129").Value = 2000
$ws.Range("J129").Value = 1431350.9
$ws.Range("K129").Value = 6000
$ws.Range("L129").Value = 4294052.699999999
$ws.Range("M129").Value = -1000
$ws.Range("N129").Value = -4304052.699999999

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2356.7144
$ws.Range("I122").Value = 2399.625
$ws.Range("J122").Value = 2299.5
$ws.Range("K122").Value = 7198.875
$ws.Range("L122").Value = 6898.5
$ws.Range("M122").Value = -4748.875
$ws.Range("N122").Value = -11798.5
$ws.Range("H126").Value = 5359.5454
$ws.Range("I126").Value = 4154
$ws.Range("J126").Value = 5811.625
$ws.Range("K126").Value = 12462
$ws.Range("L126").Value = 17434.875
$ws.Range("M126").Value = -9992
$ws.Range("N126").Value = -22374.875

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7205.875
$ws.Range("I7").Value = 4825
$ws.Range("K7").Value = 4825
$ws.Range("M7").Value = -4713
$ws.Range("H40").Value = 3605
$ws.Range("I40").Value = 2744
$ws.Range("K40").Value = 2744
$ws.Range("M40").Value = -2608
$ws.Range("H46").Value = 1367
$ws.Range("I46").Value = 1428.6923
$ws.Range("K46").Value = 1428.6923
$ws.Range("M46").Value = -1240.6923
$ws.Range("H68").Value = 2247.5
$ws.Range("I68").Value = 2195
$ws.Range("J68").Value = 2300
$ws.Range("K68").Value = 2195
$ws.Range("L68").Value = 2300
$ws.Range("M68").Value = -1446
$ws.Range("N68").Value = -3798
$ws.Range("H71").Value = 2247.5
$ws.Range("I71").Value = 2195
$ws.Range("J71").Value = 2300
$ws.Range("K71").Value = 10975
$ws.Range("L71").Value = 11500
$ws.Range("M71").Value = -7231
$ws.Range("N71").Value = -18988
$ws.Range("H122").Value = 5869.8
$ws.Range("I122").Value = 4505.619
$ws.Range("J122").Value = 7377.579
$ws.Range("K122").Value = 13516.857
$ws.Range("L122").Value = 22132.737
$ws.Range("M122").Value = -11066.857
$ws.Range("N122").Value = -27032.737
$ws.Range("H125").Value = 29500
$ws.Range("J125").Value = 29500
$ws.Range("L125").Value = 29500
$ws.Range("N125").Value = -39340
$ws.Range("H126").Value = 7205.875
$ws.Range("I126").Value = 4825
$ws.Range("K126").Value = 14475
$ws.Range("M126").Value = -12005
$ws.Range("H132").Value = 3258.0476
$ws.Range("I132").Value = 3061.5
$ws.Range("K132").Value = 9184.5
$ws.Range("M132").Value = -6654.5
$ws.Range("H136").Value = 1755.5
$ws.Range("I136").Value = 1800.5
$ws.Range("K136").Value = 5401.5
$ws.Range("M136").Value = -2851.5

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 88000
$ws.Range("I70").Value = 88000
$ws.Range("K70").Value = 88000
$ws.Range("M70").Value = -87685
$ws.Range("H73").Value = 88000
$ws.Range("I73").Value = 88000
$ws.Range("K73").Value = 88000
$ws.Range("M73").Value = -86908
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H122").Value = 3715.6956
$ws.Range("I122").Value = 2999.2
$ws.Range("K122").Value = 8997.599999999999
$ws.Range("M122").Value = -6547.599999999999
$ws.Range("H132").Value = 1980.6364
$ws.Range("I132").Value = 1809.6666
$ws.Range("K132").Value = 5428.9998
$ws.Range("M132").Value = -2898.9998
